# -----------------------------------------------------------------------
# Commit: "Thu, Jul 23, 2020  1:05:39 PM"
#
# 1) Slide 16's table (shape "Google Shape;213;p29") gets a new table
#    style applied (tableStyleId GUID change).
# 2) The two theme parts (ppt/theme/theme1.xml used by the slide master,
#    ppt/theme/theme2.xml used by the notes master) swap their color
#    schemes - theme1 becomes the "Office" palette, theme2 becomes the
#    "Integral" palette that theme1 used to have. (fontScheme/fmtScheme
#    are identical between the two themes already, so only the 12
#    clrScheme colors actually need to move.)
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- helper: hex "RRGGBB" -> the decimal BGR-ordered long that the COM
#     RGB property getter/setter uses (R + G*256 + B*65536) -------------
function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# ------------------------------------------------------------------
# 1) Table style swap on slide 16 (the table on shape 3).
# ------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shape = $slide16.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{0C177BB6-5D2F-4AFF-AB8D-634A58036B90}")
    }
}

# ------------------------------------------------------------------
# 2) Swap the slide-master theme (theme1.xml, "Integral") and the
#    notes-master theme (theme2.xml, "Office Theme") color schemes.
# ------------------------------------------------------------------
$integralColors = @(
    "000000", "FFFFFF", "455F51", "E3DED1",
    "99CB38", "63A537", "E6D024", "CC9700",
    "4EB3CF", "378DA6", "6B9F25", "B26B02"
)
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$slideMasterScheme = $p.SlideMaster.ColorScheme
$notesMasterScheme = $p.NotesMaster.ColorScheme

for ($i = 1; $i -le 12; $i++) {
    $slideMasterScheme.Colors($i).RGB = HexToRgbLong($officeColors[$i - 1])
    $notesMasterScheme.Colors($i).RGB = HexToRgbLong($integralColors[$i - 1])
}

Write-Output "done"
